$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -13.027
$ws.Range("B3").Value = 5.616
$ws.Range("C3").Value = -11.902
$ws.Range("B4").Value = 6.994
$ws.Range("D8").Value = -8.42
$ws.Range("C9").Value = -10.988
$ws.Range("A11").Value = -21.645
$ws.Range("D11").Value = -7.604000000000001
$ws.Range("A12").Value = -21.66
$ws.Range("B14").Value = 6.051
$ws.Range("D14").Value = -7.678000000000002
$ws.Range("A15").Value = -21.665
$ws.Range("C15").Value = -12.724
$ws.Range("D15").Value = -8.568999999999999
$ws.Range("D17").Value = -8.42
$ws.Range("C19").Value = -12.869
$ws.Range("C20").Value = -12.064
$ws.Range("C25").Value = -13.009
$ws.Range("B26").Value = 6.593999999999999
$ws.Range("D26").Value = -8.034000000000001
$ws.Range("A27").Value = -21.397
$ws.Range("C27").Value = -13.837
$ws.Range("A28").Value = -21.43
$ws.Range("C28").Value = -13.591
$ws.Range("C30").Value = -12.568
$ws.Range("A31").Value = -21.147
$ws.Range("B31").Value = 6.205
$ws.Range("A32").Value = -21.389
$ws.Range("C32").Value = -13.5
$ws.Range("B35").Value = 6.855
$ws.Range("A36").Value = -20.658
$ws.Range("D36").Value = -8.545
$ws.Range("B37").Value = 6.524000000000001
$ws.Range("A38").Value = -20.108
$ws.Range("B39").Value = 7.284999999999999
$ws.Range("B40").Value = 8.211
$ws.Range("D42").Value = -8.286
$ws.Range("C44").Value = -12.787
$ws.Range("B45").Value = 5.788
$ws.Range("A46").Value = -21.375
$ws.Range("C47").Value = -12.392
$ws.Range("B52").Value = 4.903
$ws.Range("A54").Value = -21.28100000000001
$ws.Range("A55").Value = -21.825
$ws.Range("A56").Value = -21.538
$ws.Range("B57").Value = 5.518
$ws.Range("C58").Value = -12.429
$ws.Range("C62").Value = -12.975
$ws.Range("D64").Value = -7.664
$ws.Range("A67").Value = -21.588
$ws.Range("D68").Value = -7.002000000000001
$ws.Range("A69").Value = -21.721
$ws.Range("A72").Value = -21.567
$ws.Range("A73").Value = -20.552
$ws.Range("C77").Value = -13.345
$ws.Range("C78").Value = -12.968
$ws.Range("D79").Value = -7.736
$ws.Range("B81").Value = 6.311
$ws.Range("A83").Value = -20.489
$ws.Range("B83").Value = 6.829000000000001
$ws.Range("C84").Value = -13.391
$ws.Range("A86").Value = -22.031
$ws.Range("C89").Value = -11.26
$ws.Range("D89").Value = -6.764999999999999
$ws.Range("A91").Value = -21.747
$ws.Range("C91").Value = -11.391
$ws.Range("C92").Value = -11.589
$ws.Range("A93").Value = -21.49
$ws.Range("C96").Value = -13.346
$ws.Range("A99").Value = -21.067
$ws.Range("B100").Value = 5.846
$ws.Range("B102").Value = 6.667
$ws.Range("C102").Value = -13.112
